$wb = $excel.ActiveWorkbook

# --- Sheet: P_valores ---
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.8864980626414685
$wsP.Range("D2").Value = 0.994489008912812
$wsP.Range("E2").Value = 0.9856185649339626
$wsP.Range("F2").Value = 0.8583128676225149

$wsP.Range("B3").Value = 0.8864980626414685
$wsP.Range("D3").Value = 0.8207769040757169
$wsP.Range("E3").Value = 0.8249384418395114
$wsP.Range("F3").Value = 0.5808894491054428

$wsP.Range("B4").Value = 0.994489008912812
$wsP.Range("C4").Value = 0.8207769040757169
$wsP.Range("E4").Value = 0.9945536183387778
$wsP.Range("F4").Value = 0.7422952320800116

$wsP.Range("B5").Value = 0.9856185649339626
$wsP.Range("C5").Value = 0.8249384418395114
$wsP.Range("D5").Value = 0.9945536183387778
$wsP.Range("F5").Value = 0.7878278934441709

$wsP.Range("B6").Value = 0.8583128676225149
$wsP.Range("C6").Value = 0.5808894491054428
$wsP.Range("D6").Value = 0.7422952320800116
$wsP.Range("E6").Value = 0.7878278934441709

# --- Sheet: Estadisticos_DM ---
$wsD = $wb.Worksheets.Item("Estadisticos_DM")

$wsD.Range("C2").Value = 0.1453608237453891
$wsD.Range("D2").Value = 0.007031406160826631
$wsD.Range("E2").Value = 0.01835003850928072
$wsD.Range("F2").Value = -0.181842590765286

$wsD.Range("B3").Value = -0.1453608237453891
$wsD.Range("D3").Value = -0.2308420285208276
$wsD.Range("E3").Value = -0.2253826025874108
$wsD.Range("F3").Value = -0.5651899753631022

$wsD.Range("B4").Value = -0.007031406160826631
$wsD.Range("C4").Value = 0.2308420285208276
$wsD.Range("E4").Value = 0.006948970353943927
$wsD.Range("F4").Value = -0.335401502127754

$wsD.Range("B5").Value = -0.01835003850928072
$wsD.Range("C5").Value = 0.2253826025874108
$wsD.Range("D5").Value = -0.006948970353943927
$wsD.Range("F5").Value = -0.2743426275790073

$wsD.Range("B6").Value = 0.181842590765286
$wsD.Range("C6").Value = 0.5651899753631022
$wsD.Range("D6").Value = 0.335401502127754
$wsD.Range("E6").Value = 0.2743426275790073

$wb.Save()
